# Ajout du tir des ennemis basiques et du regard du jeu sur le tir des autres ennemis.
# Edits the "Samuel M." worksheet: fills in row 10 with a new log entry and
# moves the active selection down to D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samuel M.")
$ws.Activate()

# New log entry for row 10.
$ws.Range("A10").Value = (Get-Date -Year 2017 -Month 12 -Day 2 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "Ajouts de contenu dans la classe des bonus. Fixatif de la classe stack. Liaison du score au HUD, mise en place des collisions des bonus. Correctif des Sujets et des Observateurs, recherche de sprites pour les Bonus et les Vaisseaux. Organisation de certaines classes, suppression de méthodes inutiles, clean-up de la mémoire (destructeurs) et codification des tests unitaires de la classe stack."
$ws.Range("D10").Value = "Le mariage du vecteur et de la stack a montré quelques mauvaises suprises qui ont été corrigées et détectée grace aux tests unitaires. J'ai aussi utilisé par accident getlocalbounds au lieu de getglobalbounds pour les collisions, ce fut déplaisant."

$ws.Rows.Item(10).RowHeight = 225

# Update the view: scroll so row 8 is at the top, select D11.
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("D11").Select()
